$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: column headers "Title" / "emails" (default/no special style)
$ws.Range("A6").Value = "Title"
$ws.Range("B6").Value = "emails"

# Row 7: first data row
$ws.Range("A7").Value = "Friday afternoon meeting"
$ws.Range("B7").Value = "bg5@mailinator.com"

# Row 5: section header "meetfromicon", merged A5:B5, centered (like row 1 "meetingbook" header)
$ws.Range("A5").Value = "meetfromicon"
$ws.Range("A5:B5").HorizontalAlignment = -4108
$ws.Range("A5:B5").Merge()

# Rows 8-10: remaining data rows
$ws.Range("A8").Value = "Afterlunch meeting"
$ws.Range("B8").Value = "bg8@mailinator.com"

$ws.Range("A9").Value = "Onsite meeting"
$ws.Range("B9").Value = "bg10@mailinator.com"

$ws.Range("A10").Value = "Offshore meeting"
$ws.Range("B10").Value = "bg12@mailinator.com"

# Format data rows (7-10) as Text, matching original data-row style
$ws.Range("A7:B10").NumberFormat = "@"

# Update selection to match the final active cell
$ws.Range("B10").Select()
